# lipidcane_spearman_1_agile.xlsx update:
# - fix pressure on sugarcane biorefinery stripping column; further work on lipidcane2g
# - Inserts a new "Feedstock consumption [ton/yr]" metric column (between "TCI [10^6*USD]"
#   and "Heat exchanger network error [%]") and refreshes the Spearman correlation table
#   with re-run values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J so the existing "Heat exchanger network error [%]" header
# (and its data) shifts from J to K, leaving J free for the new metric.
$ws.Columns("J:J").Insert()

# New header label for the inserted column (row 2).
$ws.Range("J2").Value = "Feedstock consumption [ton/yr]"

# Refresh the correlation coefficients for every parameter/metric pair (rows 4-13,
# columns C-K) with the values from the re-run analysis.
$ws.Range("C4").Value = -0.7639510987774696
$ws.Range("D4").Value = 0.9498207455186382
$ws.Range("E4").Value = -0.9780214505362635
$ws.Range("F4").Value = -0.9881992049801247
$ws.Range("G4").Value = 0.9861157360345191
$ws.Range("H4").Value = -0.9853971349283733
$ws.Range("I4").Value = 0.9383169579239482
$ws.Range("J4").Value = -0.07455936398409961
$ws.Range("K4").Value = 0.7652816320408012
$ws.Range("C5").Value = -0.06547363684092103
$ws.Range("D5").Value = 0.05502287557188931
$ws.Range("E5").Value = -0.04880822020550515
$ws.Range("F5").Value = -0.0601515037875947
$ws.Range("G5").Value = 0.07853100852049111
$ws.Range("H5").Value = -0.0647851196279907
$ws.Range("I5").Value = 0.1163894097352434
$ws.Range("J5").Value = 0.1357848946223656
$ws.Range("K5").Value = -0.0144513612840321
$ws.Range("C6").Value = 0.07710792769819247
$ws.Range("D6").Value = 0.2264696617415436
$ws.Range("E6").Value = -0.0186724668116703
$ws.Range("F6").Value = -0.06441911047776196
$ws.Range("G6").Value = 0.0607470636243027
$ws.Range("H6").Value = -0.02975624390609766
$ws.Range("I6").Value = 0.003304582614565365
$ws.Range("J6").Value = -0.06317257931448288
$ws.Range("K6").Value = 0.2507012675316883
$ws.Range("C7").Value = 0.2086117152928824
$ws.Range("D7").Value = 0.08200555013875346
$ws.Range("E7").Value = 0.247660691517288
$ws.Range("F7").Value = 0.1620685517137928
$ws.Range("G7").Value = -0.1714561668778564
$ws.Range("H7").Value = 0.2294307357683943
$ws.Range("I7").Value = 0.1273726843171079
$ws.Range("J7").Value = 0.8632430810770271
$ws.Range("K7").Value = 0.04958673966849172
$ws.Range("C8").Value = 0.3525238130953275
$ws.Range("D8").Value = 0.05234530863271583
$ws.Range("E8").Value = -0.02727218180454512
$ws.Range("F8").Value = -0.0563939098477462
$ws.Range("G8").Value = 0.04852958272556034
$ws.Range("H8").Value = -0.0376929423235581
$ws.Range("I8").Value = 0.06905272631815797
$ws.Range("J8").Value = 0.04300307507687693
$ws.Range("K8").Value = -0.0008625215630390762
$ws.Range("C9").Value = 0.4398424960624016
$ws.Range("D9").Value = -0.09131928298207456
$ws.Range("E9").Value = 0.07982449561239033
$ws.Range("F9").Value = 0.09366534163354084
$ws.Range("G9").Value = -0.09636335745143992
$ws.Range("H9").Value = 0.08894922373059329
$ws.Range("I9").Value = -0.08968424210605266
$ws.Range("J9").Value = -0.005701642541063528
$ws.Range("K9").Value = -0.007197179929498239
$ws.Range("C10").Value = 0.0112532813320333
$ws.Range("D10").Value = 0.01836795919897998
$ws.Range("E10").Value = -0.04202805070126753
$ws.Range("F10").Value = -0.0989484737118428
$ws.Range("G10").Value = 0.03083639539454768
$ws.Range("H10").Value = -0.04698267456686418
$ws.Range("I10").Value = 0.002701567539188481
$ws.Range("J10").Value = -0.01775744393609841
$ws.Range("K10").Value = 0.0456761419035476
$ws.Range("C11").Value = -0.0371814295357384
$ws.Range("D11").Value = 0.07894547363684093
$ws.Range("E11").Value = -0.03430585764644117
$ws.Range("F11").Value = -0.02637815945398635
$ws.Range("G11").Value = 0.04706384804797197
$ws.Range("H11").Value = -0.03370584264606616
$ws.Range("I11").Value = 0.1046861171529288
$ws.Range("J11").Value = 0.1039225980649516
$ws.Range("K11").Value = 0.01912397809945249
$ws.Range("C12").Value = 0.06499812495312383
$ws.Range("D12").Value = 0.07178579464486613
$ws.Range("E12").Value = 0.09856746418660468
$ws.Range("F12").Value = -0.01183529588239706
$ws.Range("G12").Value = 0.007761284448612976
$ws.Range("H12").Value = 0.02279756993924848
$ws.Range("I12").Value = -0.05253131328283208
$ws.Range("J12").Value = 0.401225530638266
$ws.Range("K12").Value = 0.003517587939698493
$ws.Range("C13").Value = -0.1751098777469437
$ws.Range("D13").Value = 0.0006960174004350109
$ws.Range("E13").Value = -0.003400585014625366
$ws.Range("F13").Value = 0.007677191929798246
$ws.Range("G13").Value = 0.01263042782029202
$ws.Range("H13").Value = -0.002031050776269407
$ws.Range("I13").Value = 0.05007125178129453
$ws.Range("J13").Value = 0.01811745293632341
$ws.Range("K13").Value = 0.008608715217880447
